$d = $word.ActiveDocument
$d.Content.Find.Execute("208÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "700÷5=", 2) | Out-Null
$d.Content.Find.Execute("522÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "969÷2=", 2) | Out-Null
$d.Content.Find.Execute("736÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "150÷9=", 2) | Out-Null
$d.Content.Find.Execute("743÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "990÷7=", 2) | Out-Null
$d.Content.Find.Execute("269÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "171÷4=", 2) | Out-Null
$d.Content.Find.Execute("570÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "942÷2=", 2) | Out-Null
$d.Content.Find.Execute("975÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "591÷6=", 2) | Out-Null
$d.Content.Find.Execute("131÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "397÷5=", 2) | Out-Null
$d.Content.Find.Execute("244÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "636÷7=", 2) | Out-Null
$d.Content.Find.Execute("539÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "664÷2=", 2) | Out-Null
$d.Content.Find.Execute("601÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "994÷8=", 2) | Out-Null
$d.Content.Find.Execute("409÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "442÷9=", 2) | Out-Null
$d.Content.Find.Execute("796÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "592÷3=", 2) | Out-Null
$d.Content.Find.Execute("878÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "381÷4=", 2) | Out-Null
$d.Content.Find.Execute("226÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "117÷7=", 2) | Out-Null
$d.Content.Find.Execute("441÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "332÷7=", 2) | Out-Null
$d.Content.Find.Execute("858÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "671÷6=", 2) | Out-Null
$d.Content.Find.Execute("540÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "690÷6=", 2) | Out-Null
$d.Content.Find.Execute("622÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "390÷5=", 2) | Out-Null
$d.Content.Find.Execute("824÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "941÷4=", 2) | Out-Null
$d.Content.Find.Execute("225÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "540÷3=", 2) | Out-Null
$d.Content.Find.Execute("928÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "666÷2=", 2) | Out-Null
$d.Content.Find.Execute("197÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "227÷7=", 2) | Out-Null
$d.Content.Find.Execute("711÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "233÷4=", 2) | Out-Null
$d.Content.Find.Execute("568÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "486÷6=", 2) | Out-Null
